$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Execute" column (C) from yes/Yes to no for these test rows
$ws.Range("C5").Value = "no"
$ws.Range("C7").Value = "no"
$ws.Range("C10").Value = "no"
$ws.Range("C12").Value = "no"
$ws.Range("C14").Value = "no"
$ws.Range("C19").Value = "no"

# Update scroll position / selection to reflect where the user was working
$ws.Range("C17").Select()
$excel.ActiveWindow.ScrollRow = 6
